$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.920.83"
$ws.Range("E2").Value = "  +6.56%  "

$ws.Range("D3").Value = "3.010.61"
$ws.Range("E3").Value = "  +3.82%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.52"
$ws.Range("E5").Value = "  +2.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.92"
$ws.Range("E6").Value = "  +13.47%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "3.006.35"
$ws.Range("E8").Value = "  +3.81%  "

$ws.Range("E9").Value = "  +3.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.93"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("E11").Value = "  +7.66%  "

$ws.Range("E12").Value = "  +6.78%  "

$ws.Range("E13").Value = "  +8.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.87"
$ws.Range("E14").Value = "  +8.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.125"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "65.937.84"
$ws.Range("E16").Value = "  +6.64%  "

$ws.Range("D17").Value = "3.511.27"
$ws.Range("E17").Value = "  +3.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("E18").Value = "  +7.28%  "

$ws.Range("D19").Value = "3.011.40"
$ws.Range("E19").Value = "  +4.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "457.88"
$ws.Range("E20").Value = "  +6.26%  "

$ws.Range("E21").Value = "  +8.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.688"
$ws.Range("E22").Value = "  +5.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.39"
$ws.Range("E23").Value = "  +7.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.53"
$ws.Range("E24").Value = "  +4.69%  "

$ws.Range("E25").Value = "  +15.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.38"
$ws.Range("E26").Value = "  +3.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.59"
$ws.Range("E27").Value = "  +5.07%  "

$ws.Range("E29").Value = "  +16.99%  "

$ws.Range("E30").Value = "  +16.18%  "

$ws.Range("E31").Value = "  +4.14%  "

$ws.Range("E32").Value = "  -7.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.05"
$ws.Range("E33").Value = "  +5.58%  "

$ws.Range("E34").Value = "  +3.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  +4.64%  "

$ws.Range("E37").Value = "  +7.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.14"
$ws.Range("E38").Value = "  +12.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.01"
$ws.Range("E39").Value = "  +6.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.84"
$ws.Range("E40").Value = "  +2.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.311"
$ws.Range("E41").Value = "  +15.97%  "

$ws.Range("E42").Value = "  +6.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.67"
$ws.Range("E43").Value = "  +8.52%  "

$ws.Range("E44").Value = "  +4.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "387.01"
$ws.Range("E45").Value = "  +11.99%  "

$ws.Range("E46").Value = "  +6.32%  "

$ws.Range("D47").Value = "2.795.55"
$ws.Range("E47").Value = "  +3.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.17"
$ws.Range("E48").Value = "  +2.72%  "

$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.03"
$ws.Range("E50").Value = "  +11.47%  "

$ws.Range("E51").Value = "  +4.11%  "
